$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Cells.Item(8,1).Value = "Volume 30   Number  10"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Weekly crime statistics table updates (rows 16-29) ---
$ws.Cells.Item(16,3).Value = 1

$ws.Cells.Item(16,4).Value = 1

$ws.Cells.Item(16,7).Value = 4

$ws.Cells.Item(16,8).Value = 75

$ws.Cells.Item(16,9).Value = 20

$ws.Cells.Item(16,10).Value = 11

$ws.Cells.Item(16,11).Value = 81.818181818181

$ws.Cells.Item(16,12).Value = 100

$ws.Cells.Item(16,13).Value = 66.666666666666

$ws.Cells.Item(16,14).Value = -81.818181818181

$ws.Cells.Item(17,3).Value = 2

$ws.Cells.Item(17,4).Value = 4

$ws.Cells.Item(17,5).Value = -50

$ws.Cells.Item(17,6).Value = 11

$ws.Cells.Item(17,7).Value = 14

$ws.Cells.Item(17,8).Value = -21.428571428571

$ws.Cells.Item(17,9).Value = 21

$ws.Cells.Item(17,10).Value = 24

$ws.Cells.Item(17,11).Value = -12.5

$ws.Cells.Item(17,12).Value = -12.5

$ws.Cells.Item(17,13).Value = 90.909090909090

$ws.Cells.Item(17,14).Value = -40

$ws.Cells.Item(18,3).Value = 2

$ws.Cells.Item(18,4).Value = 6

$ws.Cells.Item(18,6).Value = 12

$ws.Cells.Item(18,7).Value = 15

$ws.Cells.Item(18,8).Value = -20

$ws.Cells.Item(18,9).Value = 30

$ws.Cells.Item(18,10).Value = 41

$ws.Cells.Item(18,11).Value = -26.829268292682

$ws.Cells.Item(18,12).Value = 20

$ws.Cells.Item(18,13).Value = 0

$ws.Cells.Item(18,14).Value = -87.068965517241

$ws.Cells.Item(19,4).Value = 9

$ws.Cells.Item(19,5).Value = 22.222222222222

$ws.Cells.Item(19,6).Value = 50

$ws.Cells.Item(19,7).Value = 44

$ws.Cells.Item(19,8).Value = 13.636363636363

$ws.Cells.Item(19,9).Value = 117

$ws.Cells.Item(19,10).Value = 130

$ws.Cells.Item(19,11).Value = -10

$ws.Cells.Item(19,12).Value = 42.682926829268

$ws.Cells.Item(19,13).Value = -11.363636363636

$ws.Cells.Item(19,14).Value = -68.119891008174

$ws.Cells.Item(20,4).Value = 3
$ws.Cells.Item(20,4).NumberFormat = "#,##0"

$ws.Cells.Item(20,5).Value = -100
$ws.Cells.Item(20,5).NumberFormat = "#,##0.0;""-""#,##0.0"

$ws.Cells.Item(20,6).Value = 2

$ws.Cells.Item(20,7).Value = 3
$ws.Cells.Item(20,7).NumberFormat = "#,##0"

$ws.Cells.Item(20,8).Value = -33.333333333333
$ws.Cells.Item(20,8).NumberFormat = "#,##0.0;""-""#,##0.0"

$ws.Cells.Item(20,10).Value = 5

$ws.Cells.Item(20,11).Value = 40

$ws.Cells.Item(20,12).Value = -12.5

$ws.Cells.Item(20,14).Value = -94.615384615384

$ws.Cells.Item(21,3).Value = 16

$ws.Cells.Item(21,4).Value = 23

$ws.Cells.Item(21,5).Value = -30.434782608695

$ws.Cells.Item(21,7).Value = 80

$ws.Cells.Item(21,8).Value = 2.5

$ws.Cells.Item(21,9).Value = 195

$ws.Cells.Item(21,10).Value = 213

$ws.Cells.Item(21,11).Value = -8.450704225352

$ws.Cells.Item(21,12).Value = 30.872483221476

$ws.Cells.Item(21,13).Value = 4.278074866310

$ws.Cells.Item(21,14).Value = -77.765108323831

$ws.Cells.Item(22,4).Value = "'0"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(22,4).PasteSpecial(-4122)

$ws.Cells.Item(22,5).Value = "'***.*"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(22,5).PasteSpecial(-4122)

$ws.Cells.Item(22,7).Value = 2

$ws.Cells.Item(22,8).Value = 50

$ws.Cells.Item(22,9).Value = 6

$ws.Cells.Item(22,11).Value = 0

$ws.Cells.Item(22,12).Value = -14.285714285714

$ws.Cells.Item(22,13).Value = -33.333333333333

$ws.Cells.Item(24,3).Value = 18

$ws.Cells.Item(24,4).Value = 32

$ws.Cells.Item(24,5).Value = -43.75

$ws.Cells.Item(24,6).Value = 67

$ws.Cells.Item(24,7).Value = 104

$ws.Cells.Item(24,8).Value = -35.576923076923

$ws.Cells.Item(24,9).Value = 176

$ws.Cells.Item(24,10).Value = 215

$ws.Cells.Item(24,11).Value = -18.139534883720

$ws.Cells.Item(24,12).Value = -6.878306878306

$ws.Cells.Item(24,13).Value = 43.089430894308

$ws.Cells.Item(25,3).Value = 5

$ws.Cells.Item(25,4).Value = 3

$ws.Cells.Item(25,5).Value = 66.666666666666

$ws.Cells.Item(25,6).Value = 21

$ws.Cells.Item(25,7).Value = 12

$ws.Cells.Item(25,8).Value = 75

$ws.Cells.Item(25,9).Value = 48

$ws.Cells.Item(25,10).Value = 40

$ws.Cells.Item(25,11).Value = 20

$ws.Cells.Item(25,12).Value = 108.695652173913

$ws.Cells.Item(25,13).Value = 2.127659574468

$ws.Cells.Item(26,6).Value = "'0"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(26,6).PasteSpecial(-4122)

$ws.Cells.Item(27,3).Value = "'0"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(27,3).PasteSpecial(-4122)

$ws.Cells.Item(27,5).Value = -100

$ws.Cells.Item(27,6).Value = 2

$ws.Cells.Item(27,7).Value = 4

$ws.Cells.Item(27,8).Value = -50

$ws.Cells.Item(27,10).Value = 12

$ws.Cells.Item(27,11).Value = -33.333333333333

$ws.Cells.Item(27,12).Value = -42.857142857142

$ws.Cells.Item(28,4).Value = "'0"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4122)

$ws.Cells.Item(28,5).Value = "'***.*"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(28,5).PasteSpecial(-4122)

$ws.Cells.Item(29,4).Value = "'0"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(29,4).PasteSpecial(-4122)

$ws.Cells.Item(29,5).Value = "'***.*"
$ws.Cells.Item(14,1).Copy()
$ws.Cells.Item(29,5).PasteSpecial(-4122)
